$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the cells that held the "Excepção" exception-flow
# rows (rows 15-20), while keeping their existing cell formatting/styles.
$ws.Range("A15:D20").ClearContents()

# Rows 16 and 19 had an explicit 60pt row height for the wrapped exception
# text; now that the text is gone, auto-fit them back to the sheet's
# default (non-custom) row height.
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(19).AutoFit()

# Update the selected range shown when the sheet is reopened.
$ws.Range("A15:D20").Select()

$wb.Save()
